# Actualizar 02-06-2021 07-12-00
# Appends a new "Disponibilidad" check cycle (14 rows) to the end of the
# data table on Sheet1, mirroring the existing repeating 14-row block
# pattern (Odoo, Blackbox, PowerBI, Dropbox, Odoo, GEE, UtilidadesOdoo,
# Filtros Dashboard, MapStore, GeoServer, Tomcat, Shiny, Github, EZ Exporter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1080
$dateValue = 44233.2999336663

$items = @(
    @{ Name = "Odoo";              DisplayUrl = "https://www.dataintelligence-group.com/"; LinkAddress = "https://www.dataintelligence-group.com/" },
    @{ Name = "Blackbox";          DisplayUrl = "https://serviciodashboard.azurewebsites.net/"; LinkAddress = "https://serviciodashboard.azurewebsites.net/" },
    @{ Name = "PowerBI";           DisplayUrl = "https://powerbi.microsoft.com/es-es/"; LinkAddress = "https://powerbi.microsoft.com/es-es/" },
    @{ Name = "Dropbox";           DisplayUrl = "https://www.dropbox.com/"; LinkAddress = "https://www.dropbox.com/" },
    @{ Name = "Odoo";              DisplayUrl = "https://dataintelligence.store/"; LinkAddress = "https://dataintelligence.store/" },
    @{ Name = "GEE";               DisplayUrl = "https://app-data-i.users.earthengine.app/"; LinkAddress = "https://app-data-i.users.earthengine.app/" },
    @{ Name = "UtilidadesOdoo";    DisplayUrl = "https://odooutil.azurewebsites.net/"; LinkAddress = "https://odooutil.azurewebsites.net/" },
    @{ Name = "Filtros Dashboard"; DisplayUrl = "https://filtradordashboard.azurewebsites.net/"; LinkAddress = "https://filtradordashboard.azurewebsites.net/" },
    @{ Name = "MapStore";          DisplayUrl = "https://ide.dataintelligence-group.com/mapstore/#/"; LinkAddress = "https://ide.dataintelligence-group.com/mapstore/"; SubAddress = "/" },
    @{ Name = "GeoServer";         DisplayUrl = "https://ide.dataintelligence-group.com/geoserver/web/?0"; LinkAddress = "https://ide.dataintelligence-group.com/geoserver/web/?0" },
    @{ Name = "Tomcat";            DisplayUrl = "https://ide.dataintelligence-group.com/"; LinkAddress = "https://ide.dataintelligence-group.com/" },
    @{ Name = "Shiny";             DisplayUrl = "https://rpubs.com/dataintelligence/"; LinkAddress = "https://rpubs.com/dataintelligence/" },
    @{ Name = "Github";            DisplayUrl = "https://github.com/Sud-Austral/"; LinkAddress = "https://github.com/Sud-Austral/" },
    @{ Name = "EZ Exporter";       DisplayUrl = "https://ezexporter.highviewapps.com/exports/export-profile/"; LinkAddress = "https://ezexporter.highviewapps.com/exports/export-profile/" }
)

$r = $startRow
foreach ($item in $items) {
    $nameCell = $ws.Cells.Item($r, 1)
    $urlCell = $ws.Cells.Item($r, 2)
    $availCell = $ws.Cells.Item($r, 3)
    $dateCell = $ws.Cells.Item($r, 4)

    $nameCell.Value = $item.Name
    $urlCell.Value = $item.DisplayUrl
    $availCell.Value = "Disponible"

    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dateCell.Value = $dateValue

    if ($item.ContainsKey("SubAddress")) {
        $h = $ws.Hyperlinks.Add($urlCell, $item.LinkAddress, $item.SubAddress)
    } else {
        $h = $ws.Hyperlinks.Add($urlCell, $item.LinkAddress)
    }
    # Hyperlinks.Add() touches the cell's formatting; re-apply the
    # Hyperlink cell style so column B keeps the same style index
    # used by every other URL cell in the table.
    $urlCell.Style = "Hyperlink"

    $r = $r + 1
}
